# Changes of 1st april 2022
# Update the two tracking-number cells (P2, P3) on Sheet1 with new values,
# keeping them stored as text (same as the original cells) rather than
# letting them be auto-converted to numbers.
#
# We build each new value via a formula (="...") in a scratch cell, copy it,
# and paste-special the value only into the target cell - this preserves the
# destination's existing (default) cell style instead of minting a new one,
# which a direct NumberFormat="@" / quote-prefix assignment would do.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$scratch = $ws.Range("Z1")

$scratch.Formula = '="320018207767"'
$scratch.Copy()
$ws.Range("P2").PasteSpecial(-4163)

$scratch.Formula = '="320018207778"'
$scratch.Copy()
$ws.Range("P3").PasteSpecial(-4163)

$scratch.Clear()
$excel.CutCopyMode = $false
